$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Half price for "Tortila Wrap Paneer" (row 11) to 0
$ws.Range("B11").Value = 0

# Clear the Half price cell for "Quesadilla Paneer" (row 12)
$ws.Range("B12").ClearContents()

# Scroll view back to top and select B11
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B11").Select()
